$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Stephen Curry / Golden State Warriors -> Darius Garland / Cleveland Cavaliers
$ws.Range("A2").Value = "Darius Garland"
$ws.Range("B2").Value = "PG"
$ws.Range("C2").Value = "Cleveland Cavaliers"

# Row 3: Austin Reaves / Los Angeles Lakers -> Tyrese Maxey / Philadelphia 76ers
$ws.Range("A3").Value = "Tyrese Maxey"
$ws.Range("C3").Value = "Philadelphia 76ers"

# Row 5: Darius Garland / Cleveland Cavaliers -> Austin Reaves / Los Angeles Lakers
$ws.Range("A5").Value = "Austin Reaves"
$ws.Range("B5").Value = "PG,SG"
$ws.Range("C5").Value = "Los Angeles Lakers"

# Row 11: Mark Williams / Charlotte Hornets -> Jarrett Allen / Cleveland Cavaliers
$ws.Range("A11").Value = "Jarrett Allen"
$ws.Range("C11").Value = "Cleveland Cavaliers"

# Row 13: Santi Aldama / Memphis Grizzlies -> Mark Williams / Los Angeles Lakers
$ws.Range("A13").Value = "Mark Williams"
$ws.Range("B13").Value = "C"
$ws.Range("C13").Value = "Los Angeles Lakers"

# Row 14: Tyrese Maxey / Philadelphia 76ers -> Andrew Nembhard / Indiana Pacers
$ws.Range("A14").Value = "Andrew Nembhard"
$ws.Range("C14").Value = "Indiana Pacers"

# Row 15: Jarrett Allen / Cleveland Cavaliers -> Stephen Curry / Golden State Warriors
$ws.Range("A15").Value = "Stephen Curry"
$ws.Range("B15").Value = "PG,SG"
$ws.Range("C15").Value = "Golden State Warriors"
